# 2.2.3.xlsx — update the English sheet:
#  1. Retitle C1 ("2.2.3 Prevalence of anaemia ..." -> "2.2.3 Proportion of
#     women with anemia to the total population").
#  2. Add a new "2022" data column (K) mirroring the existing 2021 column
#     (J) formatting, for every data row.
#  3. Drop the stray formatted-but-empty J8 cell (header row has no J
#     column once K is introduced for data rows only).
#  4. Re-point the saved selection at M15 (matches the author's last
#     selection before publishing).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the English title in C1.
$ws.Range("C1").Value = "2.2.3 Proportion of women with anemia to the total population"

# 2. New "2022" column: copy each existing 2021 (J) cell's formatting into
#    the new K cell, then overwrite with the 2022 figure, row by row.
$ws.Range("J3").Copy($ws.Range("K3"))
$ws.Range("K3").Value = 2022

$ws.Range("J4").Copy($ws.Range("K4"))
$ws.Range("K4").Value = 1.7

$ws.Range("J5").Copy($ws.Range("K5"))

$ws.Range("J6").Copy($ws.Range("K6"))
$ws.Range("K6").Value = 2.2000000000000002

$ws.Range("J7").Copy($ws.Range("K7"))
$ws.Range("K7").Value = 1.4

# Row 8 is a section header; it no longer carries a trailing blank J cell.
$ws.Range("J8").Clear()

$ws.Range("J9").Copy($ws.Range("K9"))
$ws.Range("K9").Value = 41.1

$ws.Range("J10").Copy($ws.Range("K10"))
$ws.Range("K10").Value = 65.90209110066462

$ws.Range("J11").Copy($ws.Range("K11"))
$ws.Range("K11").Value = 55.941036331149498

$ws.Range("J12").Copy($ws.Range("K12"))
$ws.Range("K12").Value = 21.263715474839199

$ws.Range("J13").Copy($ws.Range("K13"))
$ws.Range("K13").Value = 11.351981351981353

$ws.Range("J14").Copy($ws.Range("K14"))
$ws.Range("K14").Value = 32.279274699203526

$ws.Range("J15").Copy($ws.Range("K15"))
$ws.Range("K15").Value = 36.890901250539024

$ws.Range("J16").Copy($ws.Range("K16"))
$ws.Range("K16").Value = 32.421298573536646

$ws.Range("J17").Copy($ws.Range("K17"))
$ws.Range("K17").Value = 43.227712137486577

# Last (totals) row also gets the thick bottom border the rest of row 18
# already has.
$ws.Range("J17").Copy($ws.Range("K18"))
$ws.Range("K18").Borders.Item(9).LineStyle = 1
$ws.Range("K18").Borders.Item(9).Weight = -4138
$ws.Range("K18").Value = 38.737482570668021

# 4. Move the selection, matching the saved workbook state.
$ws.Range("M15").Select()
